$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: Date | Time (7.5 hours) | Accomplished (signalR task)
$ws.Range("B25").Value = "7.5 hours"
$ws.Range("A25").Value = 44461
$ws.Range("A25").NumberFormat = $ws.Range("A24").NumberFormat
$ws.Range("C25").Value = "Task: Working on implementing signalR and sorting the data"

# Row 26: Date | Time (1 hour) | Accomplished (C# interfaces self learning)
$ws.Range("B26").Value = "1 hour"
$ws.Range("A26").Value = 44461
$ws.Range("A26").NumberFormat = $ws.Range("A24").NumberFormat
$ws.Range("C26").Value = "Self learning: Learned more about C# interfaces"

# Move the active selection to A27, matching the saved workbook state
[void]$ws.Range("A27").Select()
